$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F2
$ws.Range("F2").Value = "-"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

# Row 4
$ws.Range("B4").Value = "MCT-2A-Circuitos elétricos 2"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

# Row 6
$ws.Range("B6").Value = "MCT-2A-Circuitos elétricos 2"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "[-, -, -, 'MCT-2A-Programação de Computadores']"

# D15
$ws.Range("D15").Value = "-"

# Row 18
$ws.Range("C18").Value = "[-, -, 'ELM-1NA-Lógica de Programação', 'ELM-1NA-Acionamentos Elétricos']"
$ws.Range("E18").Value = "-"

# Row 19
$ws.Range("B19").Value = "['MEC-1NB-Comandos Eletricos', -, -, -]"
$ws.Range("C19").Value = "[-, -, 'ELM-1NA-Lógica de Programação', 'ELM-1NA-Acionamentos Elétricos']"
$ws.Range("E19").Value = "-"

# Row 20
$ws.Range("B20").Value = "['MEC-1NB-Comandos Eletricos', -, -, -]"
$ws.Range("C20").Value = "[-, -, 'ELM-1NA-Lógica de Programação', 'ELM-1NA-Acionamentos Elétricos']"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "['MEC-1NB-Comandos Eletricos', -, -, -]"
$ws.Range("C21").Value = "[-, -, 'ELM-1NA-Lógica de Programação', 'ELM-1NA-Acionamentos Elétricos']"
$ws.Range("E21").Value = "-"

$wb.Save()
